$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.037.81'
$ws.Range("E2").Value = '  +1.63%  '
$ws.Range("D3").Value = '3.741.53'
$ws.Range("E3").Value = '  +0.51%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.22'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.48%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.85'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").Value = '3.739.01'
$ws.Range("E7").Value = '  +0.54%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.536'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.06%  '
$ws.Range("E10").Value = '  -0.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.42'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.11%  '
$ws.Range("E12").Value = '  -0.83%  '
$ws.Range("E13").Value = '  -0.91%  '
$ws.Range("E14").Value = '  +0.91%  '
$ws.Range("D15").Value = '4.367.52'
$ws.Range("E15").Value = '  +0.60%  '
$ws.Range("D16").Value = '3.734.95'
$ws.Range("E16").Value = '  +0.48%  '
$ws.Range("D17").Value = '69.054.36'
$ws.Range("E17").Value = '  +1.72%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.27'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.114'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.12'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.59%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.73'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +15.76%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '492.01'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.37%  '
$ws.Range("E23").Value = '  -0.86%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000150'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.65'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.44%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.29'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.19%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.28'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.51%  '
$ws.Range("E28").Value = '  -0.23%  '
$ws.Range("E29").Value = '  -0.14%  '
$ws.Range("E30").Value = '  +1.51%  '
$ws.Range("E31").Value = '  +4.75%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.99'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.92%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.49'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.25%  '
$ws.Range("D34").Value = '3.886.61'
$ws.Range("E34").Value = '  +0.75%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.108'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.58%  '
$ws.Range("D36").Value = '3.674.43'
$ws.Range("E36").Value = '  +0.32%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("E38").Value = '  +1.11%  '
$ws.Range("E39").Value = '  -0.60%  '
$ws.Range("E40").Value = '  +0.84%  '
$ws.Range("E41").Value = '  -0.69%  '
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.94'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.93%  '
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '429.99'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.33%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '48.57'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.87%  '
$ws.Range("E45").Value = '  +0.38%  '
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '39.91'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.13%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.50'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("D50").Value = '2.776.05'
$ws.Range("E50").Value = '  +0.49%  '
$ws.Range("E51").Value = '  +0.09%  '
